# Bug fix to solar capital costs in BS calculations
# - Updates the "Capital Costs of Solar PV ($/MW)" (row 4) and
#   "Capital Costs of Solar Thermal ($/MW)" (row 10) input rows on the
#   Calculations sheet with corrected figures.
# - Clarifies the footnote text in column B for both rows.
# - Restores the sheet's view to the top-left corner with a single-cell
#   selection at A34 (previously scrolled to A7 with H47:AH47 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculations")
$wsActive = $wb.ActiveSheet

$note = "Model output, due to endogenous learning (variable Construction Cost per Unit Capacity before Construction Subsidies)"
$ws.Range("B4").Value = $note
$ws.Range("B10").Value = $note

$solarPV = @(1284590,1221260,1163880,1116630,1067190,1024270,981536,928581,881475,840697,800405,760881,750426,740565,731346,722375,713743,706003,698929,692286,686161,680812,675443,670805,666412,662174,658407,654639,651301,648071,645034,642140)
$solarThermal = @(6831840,6500520,6169200,5843880,5630240,5422800,5233450,5058220,4898610,4753980,4620220,4501030,4393600,4297490,4212530,4135500,4069500,4010230,3958590,3914720,3876450,3843370,3814240,3789350,3768370,3749030,3730660,3714590,3697200,3680150,3662120,3640820)

$col = 3
foreach ($v in $solarPV) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

$col = 3
foreach ($v in $solarThermal) {
    $ws.Cells.Item(10, $col).Value = $v
    $col = $col + 1
}

# Reset the Calculations sheet view: scroll back to the top and select A34,
# without changing which sheet is active/tabSelected in the workbook.
$ws.Activate()
$ws.Range("A34").Select()
$wsActive.Activate()
